$wb = $excel.ActiveWorkbook

# --- Update data on the "SoCDTtiNTY-psgr" sheet ---
$psgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# B2 becomes a formula: 0.076 + (0.076 - 0.0725) = 0.0795
$psgr.Range("B2").Formula = "=0.076+(0.076-0.0725)"

# D2: 0.076 -> 0.0735
$psgr.Range("D2").Value = 0.0735

# B5: 0.029 -> 0.01
$psgr.Range("B5").Value = 0.01

# E5: 0.029 -> 0.01
$psgr.Range("E5").Value = 0.01

# --- Update the active sheet / selection so the workbook re-opens on
#     "SoCDTtiNTY-psgr" with E6 selected (matching the saved view state) ---
$psgr.Activate()
$psgr.Range("E6").Select()
